$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("M2").Value = 0.06943366666666666
$ws.Range("O2").Value = 0.2790002116266049
$ws.Range("Q2").Value = 5.762138979993666
$ws.Range("R2").Value = 51.85925081994299
$ws.Range("S2").Value = 0.1252572598302913
$ws.Range("T2").Value = 0.1252572598302913
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("P3").Value = 0.7209997883733951
$ws.Range("Q3").Value = 14.89067323975233
$ws.Range("R3").Value = 134.016059157771
$ws.Range("S3").Value = 0.323693151712504
$ws.Range("T3").Value = 0.323693151712504
$ws.Range("G4").Value = 63.14058933333333
$ws.Range("I4").Value = 0.3415807409566563
$ws.Range("J4").Value = 0.3415807409566563
$ws.Range("M4").Value = 0.06943366666666666
$ws.Range("O4").Value = 0.2790002116266049
$ws.Range("Q4").Value = 4.384082632907555
$ws.Range("S4").Value = 0.09530109901447963
$ws.Range("T4").Value = 0.09530109901447963
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("P5").Value = 0.7209997883733951
$ws.Range("S5").Value = 0.2462796419421767
$ws.Range("T5").Value = 0.2462796419421767
$ws.Range("I6").Value = 0.2094688475005485
$ws.Range("J6").Value = 0.2094688475005485
$ws.Range("M6").Value = 0.06943366666666666
$ws.Range("O6").Value = 0.2790002116266049
$ws.Range("S6").Value = 0.05844185278183406
$ws.Range("T6").Value = 0.05844185278183406
$ws.Range("I7").Value = 0.2094688475005485
$ws.Range("J7").Value = 0.2094688475005485
$ws.Range("P7").Value = 0.7209997883733951
$ws.Range("Q7").Value = 6.947609539591334
$ws.Range("R7").Value = 62.52848585632201
$ws.Range("S7").Value = 0.1510269947187144
$ws.Range("T7").Value = 0.1510269947187144
